# Fix target related bugs: update branch-wise stock status figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 51; $ws.Range("C2").Value = 0; $ws.Range("D2").Value = 45; $ws.Range("E2").Value = 5; $ws.Range("F2").Value = 1; $ws.Range("G2").Value = 0; $ws.Range("H2").Value = 0
$ws.Range("B3").Value = 51; $ws.Range("C3").Value = 0; $ws.Range("D3").Value = 33; $ws.Range("E3").Value = 0; $ws.Range("F3").Value = 1; $ws.Range("G3").Value = 0; $ws.Range("H3").Value = 17
$ws.Range("B4").Value = 51; $ws.Range("C4").Value = 0; $ws.Range("D4").Value = 20; $ws.Range("E4").Value = 0; $ws.Range("F4").Value = 0; $ws.Range("G4").Value = 0; $ws.Range("H4").Value = 31
$ws.Range("B5").Value = 51; $ws.Range("C5").Value = 0; $ws.Range("D5").Value = 24; $ws.Range("E5").Value = 7; $ws.Range("F5").Value = 1; $ws.Range("G5").Value = 8; $ws.Range("H5").Value = 11
$ws.Range("B6").Value = 51; $ws.Range("C6").Value = 1; $ws.Range("D6").Value = 21; $ws.Range("E6").Value = 5; $ws.Range("F6").Value = 5; $ws.Range("G6").Value = 2; $ws.Range("H6").Value = 17
$ws.Range("B7").Value = 51; $ws.Range("C7").Value = 0; $ws.Range("D7").Value = 21; $ws.Range("E7").Value = 1; $ws.Range("F7").Value = 1; $ws.Range("G7").Value = 5; $ws.Range("H7").Value = 23
$ws.Range("B8").Value = 51; $ws.Range("C8").Value = 0; $ws.Range("D8").Value = 20; $ws.Range("E8").Value = 5; $ws.Range("F8").Value = 2; $ws.Range("G8").Value = 4; $ws.Range("H8").Value = 20
$ws.Range("B9").Value = 51; $ws.Range("C9").Value = 0; $ws.Range("D9").Value = 24; $ws.Range("E9").Value = 9; $ws.Range("F9").Value = 1; $ws.Range("G9").Value = 5; $ws.Range("H9").Value = 12
$ws.Range("B10").Value = 51; $ws.Range("C10").Value = 0; $ws.Range("D10").Value = 23; $ws.Range("E10").Value = 10; $ws.Range("F10").Value = 2; $ws.Range("G10").Value = 6; $ws.Range("H10").Value = 10
$ws.Range("B11").Value = 51; $ws.Range("C11").Value = 0; $ws.Range("D11").Value = 20; $ws.Range("E11").Value = 6; $ws.Range("F11").Value = 1; $ws.Range("G11").Value = 0; $ws.Range("H11").Value = 24
$ws.Range("B12").Value = 51; $ws.Range("C12").Value = 0; $ws.Range("D12").Value = 21; $ws.Range("E12").Value = 4; $ws.Range("F12").Value = 0; $ws.Range("G12").Value = 4; $ws.Range("H12").Value = 22
$ws.Range("B13").Value = 51; $ws.Range("C13").Value = 0; $ws.Range("D13").Value = 23; $ws.Range("E13").Value = 10; $ws.Range("F13").Value = 1; $ws.Range("G13").Value = 2; $ws.Range("H13").Value = 15
$ws.Range("B14").Value = 51; $ws.Range("C14").Value = 0; $ws.Range("D14").Value = 22; $ws.Range("E14").Value = 6; $ws.Range("F14").Value = 5; $ws.Range("G14").Value = 5; $ws.Range("H14").Value = 13
$ws.Range("B15").Value = 51; $ws.Range("C15").Value = 0; $ws.Range("D15").Value = 21; $ws.Range("E15").Value = 2; $ws.Range("F15").Value = 2; $ws.Range("G15").Value = 3; $ws.Range("H15").Value = 23
$ws.Range("B16").Value = 51; $ws.Range("C16").Value = 0; $ws.Range("D16").Value = 24; $ws.Range("E16").Value = 7; $ws.Range("F16").Value = 3; $ws.Range("G16").Value = 5; $ws.Range("H16").Value = 12
$ws.Range("B17").Value = 51; $ws.Range("C17").Value = 0; $ws.Range("D17").Value = 23; $ws.Range("E17").Value = 9; $ws.Range("F17").Value = 4; $ws.Range("G17").Value = 3; $ws.Range("H17").Value = 12
$ws.Range("B18").Value = 51; $ws.Range("C18").Value = 0; $ws.Range("D18").Value = 22; $ws.Range("E18").Value = 7; $ws.Range("F18").Value = 4; $ws.Range("G18").Value = 4; $ws.Range("H18").Value = 14
$ws.Range("B19").Value = 51; $ws.Range("C19").Value = 0; $ws.Range("D19").Value = 20; $ws.Range("E19").Value = 4; $ws.Range("F19").Value = 2; $ws.Range("G19").Value = 2; $ws.Range("H19").Value = 23
$ws.Range("B20").Value = 51; $ws.Range("C20").Value = 0; $ws.Range("D20").Value = 21; $ws.Range("E20").Value = 2; $ws.Range("F20").Value = 3; $ws.Range("G20").Value = 6; $ws.Range("H20").Value = 19
$ws.Range("B21").Value = 51; $ws.Range("C21").Value = 0; $ws.Range("D21").Value = 23; $ws.Range("E21").Value = 4; $ws.Range("F21").Value = 1; $ws.Range("G21").Value = 1; $ws.Range("H21").Value = 22
$ws.Range("B22").Value = 51; $ws.Range("C22").Value = 0; $ws.Range("D22").Value = 27; $ws.Range("E22").Value = 5; $ws.Range("F22").Value = 3; $ws.Range("G22").Value = 3; $ws.Range("H22").Value = 13
$ws.Range("B23").Value = 51; $ws.Range("C23").Value = 0; $ws.Range("D23").Value = 21; $ws.Range("E23").Value = 1; $ws.Range("F23").Value = 4; $ws.Range("G23").Value = 6; $ws.Range("H23").Value = 19
$ws.Range("B24").Value = 51; $ws.Range("C24").Value = 0; $ws.Range("D24").Value = 26; $ws.Range("E24").Value = 5; $ws.Range("F24").Value = 4; $ws.Range("G24").Value = 3; $ws.Range("H24").Value = 13
$ws.Range("B25").Value = 51; $ws.Range("C25").Value = 0; $ws.Range("D25").Value = 21; $ws.Range("E25").Value = 3; $ws.Range("F25").Value = 3; $ws.Range("G25").Value = 3; $ws.Range("H25").Value = 21
$ws.Range("B26").Value = 51; $ws.Range("C26").Value = 0; $ws.Range("D26").Value = 25; $ws.Range("E26").Value = 6; $ws.Range("F26").Value = 11; $ws.Range("G26").Value = 3; $ws.Range("H26").Value = 6
$ws.Range("B27").Value = 51; $ws.Range("C27").Value = 0; $ws.Range("D27").Value = 22; $ws.Range("E27").Value = 3; $ws.Range("F27").Value = 2; $ws.Range("G27").Value = 5; $ws.Range("H27").Value = 19
$ws.Range("B28").Value = 51; $ws.Range("C28").Value = 0; $ws.Range("D28").Value = 21; $ws.Range("E28").Value = 0; $ws.Range("F28").Value = 0; $ws.Range("G28").Value = 1; $ws.Range("H28").Value = 29
$ws.Range("B29").Value = 51; $ws.Range("C29").Value = 0; $ws.Range("D29").Value = 22; $ws.Range("E29").Value = 3; $ws.Range("F29").Value = 4; $ws.Range("G29").Value = 4; $ws.Range("H29").Value = 18
$ws.Range("B30").Value = 51; $ws.Range("C30").Value = 0; $ws.Range("D30").Value = 21; $ws.Range("E30").Value = 6; $ws.Range("F30").Value = 6; $ws.Range("G30").Value = 3; $ws.Range("H30").Value = 15
$ws.Range("B31").Value = 51; $ws.Range("C31").Value = 0; $ws.Range("D31").Value = 22; $ws.Range("E31").Value = 8; $ws.Range("F31").Value = 2; $ws.Range("G31").Value = 3; $ws.Range("H31").Value = 16
$ws.Range("B32").Value = 51; $ws.Range("C32").Value = 0; $ws.Range("D32").Value = 22; $ws.Range("E32").Value = 9; $ws.Range("F32").Value = 5; $ws.Range("G32").Value = 5; $ws.Range("H32").Value = 10
